$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.199.88"
$ws.Range("E2").Value = "  +5.88%  "
$ws.Range("D3").Value = "3.539.81"
$ws.Range("E3").Value = "  +7.36%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'189.49"
$ws.Range("E5").Value = "  +10.03%  "
$ws.Range("D6").Value = "'559.77"
$ws.Range("E6").Value = "  +5.19%  "
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.532.70"
$ws.Range("E7").Value = "  +6.79%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.618"
$ws.Range("E8").Value = "  +3.70%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'0.634"
$ws.Range("E10").Value = "  +3.74%  "
$ws.Range("E11").Value = "  +13.15%  "
$ws.Range("D12").Value = "'55.13"
$ws.Range("E12").Value = "  +3.28%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("E13").Value = "  +5.66%  "
$ws.Range("D14").Value = "'9.41"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").Value = "4.095.43"
$ws.Range("E15").Value = "  +7.43%  "
$ws.Range("D16").Value = "3.537.70"
$ws.Range("E16").Value = "  +7.67%  "
$ws.Range("E17").Value = "  +3.22%  "
$ws.Range("D18").Value = "67.101.51"
$ws.Range("E18").Value = "  +6.07%  "
$ws.Range("D19").Value = "'18.28"
$ws.Range("E19").Value = "  +4.84%  "
$ws.Range("D20").Value = "'12.12"
$ws.Range("E20").Value = "  +7.85%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  +3.46%  "
$ws.Range("D22").Value = "'431.53"
$ws.Range("E22").Value = "  +16.30%  "
$ws.Range("D23").Value = "'4.08"
$ws.Range("E23").Value = "  +8.60%  "
$ws.Range("D24").Value = "'85.57"
$ws.Range("E24").Value = "  +5.25%  "
$ws.Range("D25").Value = "'4.14"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").Value = "'11.14"
$ws.Range("E26").Value = "  -3.35%  "
$ws.Range("D27").Value = "'2.91"
$ws.Range("E27").Value = "  +8.76%  "
$ws.Range("D28").Value = "'12.43"
$ws.Range("E28").Value = "  +9.33%  "
$ws.Range("D29").Value = "'9.12"
$ws.Range("E29").Value = "  +9.99%  "
$ws.Range("D30").Value = "'30.54"
$ws.Range("E30").Value = "  +5.71%  "
$ws.Range("D31").Value = "'643.74"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'6.59"
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("D33").Value = "'11.76"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("E34").Value = "  +4.19%  "
$ws.Range("D35").Value = "'60.10"
$ws.Range("E35").Value = "  +4.92%  "
$ws.Range("D36").Value = "'38.57"
$ws.Range("E36").Value = "  +4.42%  "
$ws.Range("D37").Value = "0.0₃0816"
$ws.Range("E37").Value = "  +11.46%  "
$ws.Range("D38").Value = "'0.147"
$ws.Range("E38").Value = "  +17.79%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'0.391"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").Value = "'3.41"
$ws.Range("E41").Value = "  +16.75%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "3.067.10"
$ws.Range("E43").Value = "  +5.59%  "
$ws.Range("D44").Value = "'2.66"
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("D45").Value = "'2.89"
$ws.Range("E45").Value = "  +10.34%  "
$ws.Range("E46").Value = "  +11.77%  "
$ws.Range("D47").Value = "'0.0419"
$ws.Range("E47").Value = "  +5.38%  "
$ws.Range("E48").Value = "  +3.19%  "
$ws.Range("E49").Value = "  +4.98%  "
$ws.Range("D50").Value = "'141.85"
$ws.Range("E50").Value = "  +5.02%  "
$ws.Range("D51").Value = "'8.67"
$ws.Range("E51").Value = "  +10.74%  "
